{"js": "// Move the \"_GoBack\" bookmark from the end of the \"We split date and time...\"\n// paragraph to the end of the \"Tables were combined with an OUTER JOIN...\"\n// paragraph, and add a trailing \".\" to that paragraph (new run) right before\n// the relocated bookmark.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate paragraphs by their text content rather than hard-coded indices so\n// the script is resilient to minor structural differences.\nlet joinParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Tables\") !== -1 && t.indexOf(\"OUTER\") !== -1 && t.indexOf(\"City field\") !== -1) {\n    joinParaIndex = i;\n    break;\n  }\n}\n\nif (joinParaIndex === -1) {\n  throw new Error(\"Could not find the 'Tables ... City field' paragraph.\");\n}\n\nconst targetParagraph = paragraphs.items[joinParaIndex];\n\n// Remove the existing \"_GoBack\" bookmark (currently sitting at the end of the\n// previous paragraph).\nconst existingBookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingBookmark.load(\"isNullObject\");\nawait context.sync();\nif (!existingBookmark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Append the missing \".\" to the \"Tables ... City field\" paragraph.\nconst endRange = targetParagraph.getRange(\"End\");\nendRange.insertText(\".\", \"End\");\nawait context.sync();\n\n// Re-insert the \"_GoBack\" bookmark at the (new) end of that paragraph.\nconst newEndRange = targetParagraph.getRange(\"End\");\nnewEndRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Move the \"_GoBack\" bookmark from the end of the \"We split date and time...\"\n# paragraph to the end of the \"Tables were combined with an OUTER JOIN...\"\n# paragraph, and add a trailing \".\" to that paragraph (as a new run) right\n# before the relocated bookmark.\n\n$d = $word.ActiveDocument\n$bms = $d.Bookmarks\n\n# Locate the \"Tables ... City field\" paragraph by content instead of a\n# hard-coded index so the script is resilient to minor structural drift.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Tables*\" -and $t -like \"*OUTER*\" -and $t -like \"*City field*\") {\n        $targetParagraph = $p\n        break\n    }\n}\nif ($null -eq $targetParagraph) {\n    throw \"Could not find the 'Tables ... City field' paragraph.\"\n}\n\n# Remove the existing \"_GoBack\" bookmark (sitting at the end of the previous\n# paragraph).\nif ($bms.Exists(\"_GoBack\")) {\n    $bms.Item(\"_GoBack\").Delete()\n}\n\n# End-of-text position of the target paragraph, excluding its paragraph mark.\n$pEnd = $targetParagraph.Range.End - 1\n\n# Append the missing \".\" as its own run. Inserting directly at the paragraph's\n# true end (i.e. right before the paragraph mark, with nothing else following\n# in this paragraph) keeps it a separate <w:r> instead of merging into the\n# preceding run.\n$periodRange = $d.Range($pEnd, $pEnd)\n$periodRange.InsertAfter(\".\")\n\n# Re-resolve the (now one character later) true paragraph end and insert a\n# temporary one-character placeholder there. This is a workaround for a COM\n# bridge quirk: adding a bookmark via a zero-length Range whose Start/End sit\n# exactly AT a paragraph's true end position gets mis-resolved to the start of\n# the document. Keeping one throwaway character after the insertion point\n# means the bookmark target is no longer that special boundary position.\n$pEnd2 = $targetParagraph.Range.End - 1\n$placeholderRange = $d.Range($pEnd2, $pEnd2)\n$placeholderRange.InsertAfter(\"X\")\n\n# Insert the relocated bookmark right between \".\" and the placeholder \"X\".\n$bookmarkTarget = $d.Range($pEnd2, $pEnd2)\n$bms.Add(\"_GoBack\", $bookmarkTarget) | Out-Null\n\n# Remove the placeholder character.\n$placeholderCleanup = $d.Range($pEnd2, $pEnd2 + 1)\n$placeholderCleanup.Delete()\n"}
